$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("175").Insert()

$ws.Range("A175").Value = 3
$ws.Range("B175").Value = "Femacal de La Calera"
$ws.Range("C175").Value = "Coquimbo"
$ws.Range("D175").Value = 44603
$ws.Range("E175").Value = 5
$ws.Range("F175").Value = 100112009
$ws.Range("G175").Value = "Acelga"
$ws.Range("H175").Value = "Sin especificar"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 230
$ws.Range("K175").Value = 2300
$ws.Range("L175").Value = 2500
$ws.Range("M175").Value = 2396
$ws.Range("N175").Value = "$/docena de atados (6 kilos)"
$ws.Range("O175").Value = "Provincia de Quillota"
$ws.Range("P175").Value = 399
$ws.Range("Q175").Value = 6
$ws.Range("R175").Value = "Hortaliza"
